$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.164.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.48%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.911.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.04%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9969"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.33%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7423"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.42%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9972"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3124"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.66%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.71"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.95%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06962"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.63%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7806"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07965"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.874.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.80%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.269"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.02%  "

$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.04%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.166.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.860"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007830"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9972"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.145.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.85%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9964"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.36%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.996"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.94%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.400"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.48"
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1283"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.069"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.340"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.72%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.539"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.75%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.334"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.099"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05140"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.298"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7422"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.54%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.711"
$ws.Range("D37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01943"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.797"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.333"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4490"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.954"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.50%  "

$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.837"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.00%  "

$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9982"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.21%  "

$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.8380"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.85%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.943"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.17%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.44%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.56%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.051.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.80%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "938.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.34%  "
